# Apply the CodeSystem-ValidAgeReason.xlsx metadata update:
#  - URL value: pythia -> cicada
#  - Date value refreshed
#  - New "Jurisdiction" row inserted (empty value) right after "Contact" / before "Description"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new row above row 11 ("Description"), pushing everything down by one.
$ws.Rows.Item(11).Insert()

# Fill the newly inserted row 11 with the Jurisdiction property (value left blank)
# and match the formatting used by the rest of the property rows (e.g. row 12)
# by copying just the cell formats (not a full row copy, which would widen the
# used range to the whole sheet).
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""

# Update the URL value (row 2, column B).
$ws.Cells.Item(2, 2).Value = "http://fhirfli.dev/fhir/ig/cicada/CodeSystem/ValidAgeReason"

# Update the Date value (row 8, column B).
$ws.Cells.Item(8, 2).Value = "2026-02-11T14:37:07-05:00"
